$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update existing rows 5-42 (shuffled colors/durations/sounds/targets)
$ws.Range("C5").Value = 0.05
$ws.Range("F5").Value = "silent.wav"
$ws.Range("B6").Value = "orange.png"
$ws.Range("F6").Value = "beep.wav"
$ws.Range("B9").Value = "green.png"
$ws.Range("F9").Value = "silent.wav"
$ws.Range("B10").Value = "orange.png"
$ws.Range("F10").Value = "beep.wav"
$ws.Range("B11").Value = "green.png"
$ws.Range("F11").Value = "silent.wav"
$ws.Range("B12").Value = "orange.png"
$ws.Range("F12").Value = "beep.wav"
$ws.Range("B13").Value = "green.png"
$ws.Range("D13").Value = 0
$ws.Range("E13").ClearContents()
$ws.Range("F13").Value = "silent.wav"
$ws.Range("B14").Value = "red.png"
$ws.Range("D14").Value = 1
$ws.Range("E14").Value = 1
$ws.Range("F14").Value = "beep.wav"
$ws.Range("B17").Value = "green.png"
$ws.Range("F17").Value = "silent.wav"
$ws.Range("B18").Value = "orange.png"
$ws.Range("F18").Value = "beep.wav"
$ws.Range("B19").Value = "green.png"
$ws.Range("D19").Value = 0
$ws.Range("E19").ClearContents()
$ws.Range("F19").Value = "silent.wav"
$ws.Range("B20").Value = "red.png"
$ws.Range("D20").Value = 1
$ws.Range("E20").Value = 1
$ws.Range("B21").Value = "orange.png"
$ws.Range("C21").Value = 0.05
$ws.Range("F21").Value = "silent.wav"
$ws.Range("B22").Value = "orange.png"
$ws.Range("F22").Value = "beep.wav"
$ws.Range("B23").Value = "red.png"
$ws.Range("C23").Value = 0.05
$ws.Range("F23").Value = "silent.wav"
$ws.Range("D24").Value = 0
$ws.Range("E24").ClearContents()
$ws.Range("B26").Value = "orange.png"
$ws.Range("F26").Value = "beep.wav"
$ws.Range("B27").Value = "red.png"
$ws.Range("C27").Value = 0.05
$ws.Range("E27").Value = 0
$ws.Range("B28").Value = "red.png"
$ws.Range("D28").Value = 1
$ws.Range("E28").Value = 1
$ws.Range("B30").Value = "green.png"
$ws.Range("D30").Value = 0
$ws.Range("E30").ClearContents()
$ws.Range("F30").Value = "silent.wav"
$ws.Range("B31").Value = "green.png"
$ws.Range("F31").Value = "silent.wav"
$ws.Range("B32").Value = "orange.png"
$ws.Range("B34").Value = "red.png"
$ws.Range("D34").Value = 1
$ws.Range("E34").Value = 1
$ws.Range("F34").Value = "beep.wav"
$ws.Range("C35").Value = 0.05
$ws.Range("F35").Value = "silent.wav"
$ws.Range("B36").Value = "orange.png"
$ws.Range("D36").Value = 0
$ws.Range("E36").ClearContents()
$ws.Range("B37").Value = "red.png"
$ws.Range("C37").Value = 0.05
$ws.Range("D38").Value = 0
$ws.Range("E38").ClearContents()
$ws.Range("B40").Value = "green.png"
$ws.Range("D40").Value = 0
$ws.Range("E40").ClearContents()
$ws.Range("F40").Value = "silent.wav"
$ws.Range("B41").Value = "orange.png"
$ws.Range("F41").Value = "beep.wav"
$ws.Range("B42").Value = "red.png"
$ws.Range("C42").Value = 0.05
$ws.Range("E42").Value = 0

# New rows 43-49 appended to the bottom of the table
$ws.Range("A43").Value = 1
$ws.Range("B43").Value = "red.png"
$ws.Range("C43").Value = 3
$ws.Range("D43").Value = 1
$ws.Range("E43").Value = 1
$ws.Range("F43").Value = "beep.wav"
$ws.Range("A44").Value = 1
$ws.Range("B44").Value = "green.png"
$ws.Range("C44").Value = 3
$ws.Range("D44").Value = 0
$ws.Range("F44").Value = "silent.wav"
$ws.Range("A45").Value = 1
$ws.Range("B45").Value = "red.png"
$ws.Range("C45").Value = 3
$ws.Range("D45").Value = 1
$ws.Range("E45").Value = 1
$ws.Range("F45").Value = "beep.wav"
$ws.Range("A46").Value = 1
$ws.Range("B46").Value = "green.png"
$ws.Range("C46").Value = 3
$ws.Range("D46").Value = 0
$ws.Range("F46").Value = "silent.wav"
$ws.Range("A47").Value = 1
$ws.Range("B47").Value = "red.png"
$ws.Range("C47").Value = 3
$ws.Range("D47").Value = 1
$ws.Range("E47").Value = 1
$ws.Range("F47").Value = "beep.wav"
$ws.Range("A48").Value = 1
$ws.Range("B48").Value = "green.png"
$ws.Range("C48").Value = 3
$ws.Range("D48").Value = 0
$ws.Range("F48").Value = "silent.wav"
$ws.Range("A49").Value = 1
$ws.Range("B49").Value = "green.png"
$ws.Range("C49").Value = 2.9
$ws.Range("D49").Value = 0
$ws.Range("F49").Value = "silent.wav"

# Leave the selection on C5 (matches the saved sheet view)
$ws.Range("C5").Select()
